# repull data, push all data, mean calculation
# Update column F ("dSF") values on Sheet1 for the rows that changed
# after re-pulling the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = 5
    3  = 1
    4  = 2
    6  = -3
    7  = -3
    8  = -2
    9  = -1
    10 = -1
    11 = -1
    12 = -3
    13 = -2
    14 = -4
    15 = -2
    16 = -1
    17 = -2
    18 = 2
    19 = -1
    22 = -2
    25 = 1
    26 = -5
    27 = -1
    28 = -3
    29 = -3
    30 = -1
    31 = -1
    32 = 2
    33 = -11
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
